# tdf#144092 test-fixture touch-up:
#   1. Refresh the cached "datetimeFigureOut" date placeholder text (slide
#      master + every slide layout) from 8/30/2022 to 10/20/2023.
#   2. Swap the theme's major/minor Latin fonts from Calibri Light/Calibri
#      to Liberation Sans (a font bundled with LibreOffice) so the test is
#      robust across systems that don't have Calibri installed.

$p = $ppt.ActivePresentation

function Update-DateShapes($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "8/30/2022") {
                $shp.TextFrame.TextRange.Text = "10/20/2023"
            }
        }
    }
}

# Slide master's own date placeholder.
Update-DateShapes $p.SlideMaster.Shapes

# Every slide layout has its own copy of the date placeholder too.
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    Update-DateShapes $layout.Shapes
}

# Theme fonts: Calibri Light / Calibri -> Liberation Sans.
$fontScheme = $p.SlideMaster.Theme.ThemeElements.ThemeFontScheme
$fontScheme.MajorFont.Latin = "Liberation Sans"
$fontScheme.MinorFont.Latin = "Liberation Sans"
